$wb = $excel.ActiveWorkbook

# ----- Sheet: LP1912 -----
$ws = $wb.Worksheets.Item('LP1912')
$ws.Range("A2").Value = 'Última actualización: 10:26:41'
$ws.Range("A3").Value = 'Total filas: 173'
$arr = New-Object 'object[,]' 173,5
$arr[0,0] = '04:18:53'
$arr[0,1] = '04:45'
$arr[0,2] = '215A_EL PATO'
$arr[0,3] = 27
$arr[0,4] = 'LP1912'
$arr[1,0] = '04:40:33'
$arr[1,1] = '04:46'
$arr[1,2] = '215A_EL PATO'
$arr[1,3] = 6
$arr[1,4] = 'LP1912'
$arr[2,0] = '04:18:53'
$arr[2,1] = '04:53'
$arr[2,2] = '11_ETCHEVERRY'
$arr[2,3] = 35
$arr[2,4] = 'LP1912'
$arr[3,0] = '04:54:03'
$arr[3,1] = '04:55'
$arr[3,2] = '11_ETCHEVERRY'
$arr[3,3] = 1
$arr[3,4] = 'LP1912'
$arr[4,0] = '04:54:03'
$arr[4,1] = '05:15'
$arr[4,2] = '14_ABASTO'
$arr[4,3] = 21
$arr[4,4] = 'LP1912'
$arr[5,0] = '04:18:53'
$arr[5,1] = '05:16'
$arr[5,2] = '17_ROMERO'
$arr[5,3] = 58
$arr[5,4] = 'LP1912'
$arr[6,0] = '04:40:33'
$arr[6,1] = '05:16'
$arr[6,2] = '14_ABASTO'
$arr[6,3] = 36
$arr[6,4] = 'LP1912'
$arr[7,0] = '05:18:56'
$arr[7,1] = '05:20'
$arr[7,2] = '14_ABASTO'
$arr[7,3] = 2
$arr[7,4] = 'LP1912'
$arr[8,0] = '04:18:53'
$arr[8,1] = '05:21'
$arr[8,2] = '23_HERNANDEZ'
$arr[8,3] = 63
$arr[8,4] = 'LP1912'
$arr[9,0] = '04:40:33'
$arr[9,1] = '05:22'
$arr[9,2] = '23_HERNANDEZ'
$arr[9,3] = 42
$arr[9,4] = 'LP1912'
$arr[10,0] = '04:18:53'
$arr[10,1] = '05:34'
$arr[10,2] = '215B_EL PATO'
$arr[10,3] = 76
$arr[10,4] = 'LP1912'
$arr[11,0] = '04:54:03'
$arr[11,1] = '05:35'
$arr[11,2] = '215B_EL PATO'
$arr[11,3] = 41
$arr[11,4] = 'LP1912'
$arr[12,0] = '04:18:53'
$arr[12,1] = '05:37'
$arr[12,2] = '14_ABASTO'
$arr[12,3] = 79
$arr[12,4] = 'LP1912'
$arr[13,0] = '04:18:53'
$arr[13,1] = '05:46'
$arr[13,2] = '15_ABASTO'
$arr[13,3] = 88
$arr[13,4] = 'LP1912'
$arr[14,0] = '05:49:40'
$arr[14,1] = '05:51'
$arr[14,2] = '14_ABASTO'
$arr[14,3] = 2
$arr[14,4] = 'LP1912'
$arr[15,0] = '04:40:33'
$arr[15,1] = '06:04'
$arr[15,2] = '16_SANTA ANA'
$arr[15,3] = 84
$arr[15,4] = 'LP1912'
$arr[16,0] = '05:49:40'
$arr[16,1] = '06:06'
$arr[16,2] = '16_SANTA ANA'
$arr[16,3] = 17
$arr[16,4] = 'LP1912'
$arr[17,0] = '04:18:53'
$arr[17,1] = '06:07'
$arr[17,2] = '16_SANTA ANA'
$arr[17,3] = 109
$arr[17,4] = 'LP1912'
$arr[18,0] = '04:18:53'
$arr[18,1] = '06:11'
$arr[18,2] = '215A_EL PATO'
$arr[18,3] = 113
$arr[18,4] = 'LP1912'
$arr[19,0] = '04:54:03'
$arr[19,1] = '06:12'
$arr[19,2] = '215A_EL PATO'
$arr[19,3] = 78
$arr[19,4] = 'LP1912'
$arr[20,0] = '04:18:53'
$arr[20,1] = '06:13'
$arr[20,2] = '225_HARAS DEL SUR'
$arr[20,3] = 115
$arr[20,4] = 'LP1912'
$arr[21,0] = '04:40:33'
$arr[21,1] = '06:14'
$arr[21,2] = '225_HARAS DEL SUR'
$arr[21,3] = 94
$arr[21,4] = 'LP1912'
$arr[22,0] = '06:15:04'
$arr[22,1] = '06:15'
$arr[22,2] = '225_HARAS DEL SUR'
$arr[22,3] = 0
$arr[22,4] = 'LP1912'
$arr[23,0] = '05:18:56'
$arr[23,1] = '06:20'
$arr[23,2] = '26_HERNANDEZ'
$arr[23,3] = 62
$arr[23,4] = 'LP1912'
$arr[24,0] = '04:40:33'
$arr[24,1] = '06:21'
$arr[24,2] = '26_HERNANDEZ'
$arr[24,3] = 101
$arr[24,4] = 'LP1912'
$arr[25,0] = '05:18:56'
$arr[25,1] = '06:26'
$arr[25,2] = '23_HERNANDEZ'
$arr[25,3] = 68
$arr[25,4] = 'LP1912'
$arr[26,0] = '04:40:33'
$arr[26,1] = '06:27'
$arr[26,2] = '23_HERNANDEZ'
$arr[26,3] = 107
$arr[26,4] = 'LP1912'
$arr[27,0] = '06:15:04'
$arr[27,1] = '06:28'
$arr[27,2] = '23_HERNANDEZ'
$arr[27,3] = 13
$arr[27,4] = 'LP1912'
$arr[28,0] = '04:40:33'
$arr[28,1] = '06:29'
$arr[28,2] = '86_EST CHICA-ESC AGRARIA'
$arr[28,3] = 109
$arr[28,4] = 'LP1912'
$arr[29,0] = '04:54:03'
$arr[29,1] = '06:30'
$arr[29,2] = '86_EST CHICA-ESC AGRARIA'
$arr[29,3] = 96
$arr[29,4] = 'LP1912'
$arr[30,0] = '04:40:33'
$arr[30,1] = '06:31'
$arr[30,2] = '16_SANTA ANA'
$arr[30,3] = 111
$arr[30,4] = 'LP1912'
$arr[31,0] = '05:18:56'
$arr[31,1] = '06:43'
$arr[31,2] = '225_C ROCA-H SUR'
$arr[31,3] = 85
$arr[31,4] = 'LP1912'
$arr[32,0] = '04:54:03'
$arr[32,1] = '06:44'
$arr[32,2] = '225_C ROCA-H SUR'
$arr[32,3] = 110
$arr[32,4] = 'LP1912'
$arr[33,0] = '05:18:56'
$arr[33,1] = '06:46'
$arr[33,2] = '215C_EL PATO'
$arr[33,3] = 88
$arr[33,4] = 'LP1912'
$arr[34,0] = '06:43:40'
$arr[34,1] = '06:46'
$arr[34,2] = '225_C ROCA-H SUR'
$arr[34,3] = 3
$arr[34,4] = 'LP1912'
$arr[35,0] = '04:54:03'
$arr[35,1] = '06:47'
$arr[35,2] = '215C_EL PATO'
$arr[35,3] = 113
$arr[35,4] = 'LP1912'
$arr[36,0] = '06:57:30'
$arr[36,1] = '06:58'
$arr[36,2] = '14_ABASTO'
$arr[36,3] = 1
$arr[36,4] = 'LP1912'
$arr[37,0] = '05:18:56'
$arr[37,1] = '06:58'
$arr[37,2] = '10_OLMOS'
$arr[37,3] = 100
$arr[37,4] = 'LP1912'
$arr[38,0] = '05:18:56'
$arr[38,1] = '06:59'
$arr[38,2] = '14_ABASTO'
$arr[38,3] = 101
$arr[38,4] = 'LP1912'
$arr[39,0] = '06:15:04'
$arr[39,1] = '07:00'
$arr[39,2] = '14_ABASTO'
$arr[39,3] = 45
$arr[39,4] = 'LP1912'
$arr[40,0] = '06:15:04'
$arr[40,1] = '07:01'
$arr[40,2] = '16_SANTA ANA'
$arr[40,3] = 46
$arr[40,4] = 'LP1912'
$arr[41,0] = '05:18:56'
$arr[41,1] = '07:04'
$arr[41,2] = '15_ABASTO'
$arr[41,3] = 106
$arr[41,4] = 'LP1912'
$arr[42,0] = '05:49:40'
$arr[42,1] = '07:04'
$arr[42,2] = '23_HERNANDEZ'
$arr[42,3] = 75
$arr[42,4] = 'LP1912'
$arr[43,0] = '05:49:40'
$arr[43,1] = '07:05'
$arr[43,2] = '15_ABASTO'
$arr[43,3] = 76
$arr[43,4] = 'LP1912'
$arr[44,0] = '05:18:56'
$arr[44,1] = '07:06'
$arr[44,2] = '225_GOMEZ'
$arr[44,3] = 108
$arr[44,4] = 'LP1912'
$arr[45,0] = '05:49:40'
$arr[45,1] = '07:07'
$arr[45,2] = '225_GOMEZ'
$arr[45,3] = 78
$arr[45,4] = 'LP1912'
$arr[46,0] = '05:18:56'
$arr[46,1] = '07:11'
$arr[46,2] = '215A_EL PATO'
$arr[46,3] = 113
$arr[46,4] = 'LP1912'
$arr[47,0] = '06:15:04'
$arr[47,1] = '07:12'
$arr[47,2] = '215A_EL PATO'
$arr[47,3] = 57
$arr[47,4] = 'LP1912'
$arr[48,0] = '05:18:56'
$arr[48,1] = '07:15'
$arr[48,2] = '11_ETCHEVERRY'
$arr[48,3] = 117
$arr[48,4] = 'LP1912'
$arr[49,0] = '06:43:40'
$arr[49,1] = '07:16'
$arr[49,2] = '16_SANTA ANA'
$arr[49,3] = 33
$arr[49,4] = 'LP1912'
$arr[50,0] = '06:15:04'
$arr[50,1] = '07:16'
$arr[50,2] = '11_ETCHEVERRY'
$arr[50,3] = 61
$arr[50,4] = 'LP1912'
$arr[51,0] = '07:20:40'
$arr[51,1] = '07:20'
$arr[51,2] = '10_OLMOS'
$arr[51,3] = 0
$arr[51,4] = 'LP1912'
$arr[52,0] = '05:49:40'
$arr[52,1] = '07:21'
$arr[52,2] = '26_HERNANDEZ'
$arr[52,3] = 92
$arr[52,4] = 'LP1912'
$arr[53,0] = '06:15:04'
$arr[53,1] = '07:23'
$arr[53,2] = '10_OLMOS'
$arr[53,3] = 68
$arr[53,4] = 'LP1912'
$arr[54,0] = '05:49:40'
$arr[54,1] = '07:29'
$arr[54,2] = '10_OLMOS'
$arr[54,3] = 100
$arr[54,4] = 'LP1912'
$arr[55,0] = '05:49:40'
$arr[55,1] = '07:31'
$arr[55,2] = '11_ETCHEVERRY'
$arr[55,3] = 102
$arr[55,4] = 'LP1912'
$arr[56,0] = '06:15:04'
$arr[56,1] = '07:32'
$arr[56,2] = '11_ETCHEVERRY'
$arr[56,3] = 77
$arr[56,4] = 'LP1912'
$arr[57,0] = '05:49:40'
$arr[57,1] = '07:32'
$arr[57,2] = '84_COLONIA URQUIZA-ESC 49'
$arr[57,3] = 103
$arr[57,4] = 'LP1912'
$arr[58,0] = '07:20:40'
$arr[58,1] = '07:34'
$arr[58,2] = '23_HERNANDEZ'
$arr[58,3] = 14
$arr[58,4] = 'LP1912'
$arr[59,0] = '05:49:40'
$arr[59,1] = '07:36'
$arr[59,2] = '27_EL RETIRO'
$arr[59,3] = 107
$arr[59,4] = 'LP1912'
$arr[60,0] = '06:15:04'
$arr[60,1] = '07:37'
$arr[60,2] = '27_EL RETIRO'
$arr[60,3] = 82
$arr[60,4] = 'LP1912'
$arr[61,0] = '05:49:40'
$arr[61,1] = '07:39'
$arr[61,2] = '10_OLMOS'
$arr[61,3] = 110
$arr[61,4] = 'LP1912'
$arr[62,0] = '07:20:40'
$arr[62,1] = '07:46'
$arr[62,2] = '16_SANTA ANA'
$arr[62,3] = 26
$arr[62,4] = 'LP1912'
$arr[63,0] = '06:43:40'
$arr[63,1] = '07:47'
$arr[63,2] = '14_ABASTO'
$arr[63,3] = 64
$arr[63,4] = 'LP1912'
$arr[64,0] = '06:15:04'
$arr[64,1] = '07:48'
$arr[64,2] = '14_ABASTO'
$arr[64,3] = 93
$arr[64,4] = 'LP1912'
$arr[65,0] = '06:43:40'
$arr[65,1] = '07:51'
$arr[65,2] = '215D_EL PATO'
$arr[65,3] = 68
$arr[65,4] = 'LP1912'
$arr[66,0] = '06:15:04'
$arr[66,1] = '07:52'
$arr[66,2] = '215D_EL PATO'
$arr[66,3] = 97
$arr[66,4] = 'LP1912'
$arr[67,0] = '07:47:32'
$arr[67,1] = '07:55'
$arr[67,2] = '10_OLMOS'
$arr[67,3] = 8
$arr[67,4] = 'LP1912'
$arr[68,0] = '07:20:40'
$arr[68,1] = '07:58'
$arr[68,2] = '16_SANTA ANA'
$arr[68,3] = 38
$arr[68,4] = 'LP1912'
$arr[69,0] = '07:20:40'
$arr[69,1] = '07:59'
$arr[69,2] = '23_HERNANDEZ'
$arr[69,3] = 39
$arr[69,4] = 'LP1912'
$arr[70,0] = '06:15:04'
$arr[70,1] = '08:01'
$arr[70,2] = '23_HERNANDEZ'
$arr[70,3] = 106
$arr[70,4] = 'LP1912'
$arr[71,0] = '06:43:40'
$arr[71,1] = '08:03'
$arr[71,2] = '23_HERNANDEZ'
$arr[71,3] = 80
$arr[71,4] = 'LP1912'
$arr[72,0] = '07:20:40'
$arr[72,1] = '08:03'
$arr[72,2] = '11_ETCHEVERRY'
$arr[72,3] = 43
$arr[72,4] = 'LP1912'
$arr[73,0] = '06:57:30'
$arr[73,1] = '08:06'
$arr[73,2] = '23_HERNANDEZ'
$arr[73,3] = 69
$arr[73,4] = 'LP1912'
$arr[74,0] = '07:47:32'
$arr[74,1] = '08:10'
$arr[74,2] = '16_SANTA ANA'
$arr[74,3] = 23
$arr[74,4] = 'LP1912'
$arr[75,0] = '06:15:04'
$arr[75,1] = '08:12'
$arr[75,2] = '15_ABASTO'
$arr[75,3] = 117
$arr[75,4] = 'LP1912'
$arr[76,0] = '07:47:32'
$arr[76,1] = '08:13'
$arr[76,2] = '10_OLMOS'
$arr[76,3] = 26
$arr[76,4] = 'LP1912'
$arr[77,0] = '07:47:32'
$arr[77,1] = '08:16'
$arr[77,2] = '26_HERNANDEZ'
$arr[77,3] = 29
$arr[77,4] = 'LP1912'
$arr[78,0] = '06:43:40'
$arr[78,1] = '08:21'
$arr[78,2] = '26_HERNANDEZ'
$arr[78,3] = 98
$arr[78,4] = 'LP1912'
$arr[79,0] = '06:43:40'
$arr[79,1] = '08:22'
$arr[79,2] = '16_P MOR-SANTA ANA'
$arr[79,3] = 99
$arr[79,4] = 'LP1912'
$arr[80,0] = '06:43:40'
$arr[80,1] = '08:23'
$arr[80,2] = '215B_EL PATO'
$arr[80,3] = 100
$arr[80,4] = 'LP1912'
$arr[81,0] = '06:43:40'
$arr[81,1] = '08:27'
$arr[81,2] = '84_COLONIA URQUIZA-ESC 49'
$arr[81,3] = 104
$arr[81,4] = 'LP1912'
$arr[82,0] = '07:47:32'
$arr[82,1] = '08:31'
$arr[82,2] = '23_HERNANDEZ'
$arr[82,3] = 44
$arr[82,4] = 'LP1912'
$arr[83,0] = '07:59:28'
$arr[83,1] = '08:33'
$arr[83,2] = '10_OLMOS'
$arr[83,3] = 34
$arr[83,4] = 'LP1912'
$arr[84,0] = '07:59:28'
$arr[84,1] = '08:34'
$arr[84,2] = '23_HERNANDEZ'
$arr[84,3] = 35
$arr[84,4] = 'LP1912'
$arr[85,0] = '07:59:28'
$arr[85,1] = '08:39'
$arr[85,2] = '26_HERNANDEZ'
$arr[85,3] = 40
$arr[85,4] = 'LP1912'
$arr[86,0] = '06:43:40'
$arr[86,1] = '08:42'
$arr[86,2] = '81_EL PELIGRO'
$arr[86,3] = 119
$arr[86,4] = 'LP1912'
$arr[87,0] = '07:20:40'
$arr[87,1] = '08:43'
$arr[87,2] = '14_ABASTO'
$arr[87,3] = 83
$arr[87,4] = 'LP1912'
$arr[88,0] = '06:57:30'
$arr[88,1] = '08:54'
$arr[88,2] = '17_ROMERO'
$arr[88,3] = 117
$arr[88,4] = 'LP1912'
$arr[89,0] = '08:57:13'
$arr[89,1] = '08:59'
$arr[89,2] = '11_ETCHEVERRY'
$arr[89,3] = 2
$arr[89,4] = 'LP1912'
$arr[90,0] = '07:20:40'
$arr[90,1] = '09:01'
$arr[90,2] = '215A_EL PATO'
$arr[90,3] = 101
$arr[90,4] = 'LP1912'
$arr[91,0] = '08:21:50'
$arr[91,1] = '09:01'
$arr[91,2] = '23_HERNANDEZ'
$arr[91,3] = 40
$arr[91,4] = 'LP1912'
$arr[92,0] = '08:57:13'
$arr[92,1] = '09:02'
$arr[92,2] = '215A_EL PATO'
$arr[92,3] = 5
$arr[92,4] = 'LP1912'
$arr[93,0] = '07:59:28'
$arr[93,1] = '09:03'
$arr[93,2] = '11_ETCHEVERRY'
$arr[93,3] = 64
$arr[93,4] = 'LP1912'
$arr[94,0] = '08:39:44'
$arr[94,1] = '09:04'
$arr[94,2] = '23_HERNANDEZ'
$arr[94,3] = 25
$arr[94,4] = 'LP1912'
$arr[95,0] = '08:57:13'
$arr[95,1] = '09:05'
$arr[95,2] = '23_HERNANDEZ'
$arr[95,3] = 8
$arr[95,4] = 'LP1912'
$arr[96,0] = '08:21:50'
$arr[96,1] = '09:07'
$arr[96,2] = '26_HERNANDEZ'
$arr[96,3] = 46
$arr[96,4] = 'LP1912'
$arr[97,0] = '07:20:40'
$arr[97,1] = '09:10'
$arr[97,2] = '16_P MOR-SANTA ANA'
$arr[97,3] = 110
$arr[97,4] = 'LP1912'
$arr[98,0] = '08:57:13'
$arr[98,1] = '09:11'
$arr[98,2] = '16_P MOR-SANTA ANA'
$arr[98,3] = 14
$arr[98,4] = 'LP1912'
$arr[99,0] = '08:21:50'
$arr[99,1] = '09:13'
$arr[99,2] = '10_OLMOS'
$arr[99,3] = 52
$arr[99,4] = 'LP1912'
$arr[100,0] = '07:20:40'
$arr[100,1] = '09:16'
$arr[100,2] = '27_EL RETIRO'
$arr[100,3] = 116
$arr[100,4] = 'LP1912'
$arr[101,0] = '08:57:13'
$arr[101,1] = '09:17'
$arr[101,2] = '27_EL RETIRO'
$arr[101,3] = 20
$arr[101,4] = 'LP1912'
$arr[102,0] = '08:21:50'
$arr[102,1] = '09:21'
$arr[102,2] = '26_HERNANDEZ'
$arr[102,3] = 60
$arr[102,4] = 'LP1912'
$arr[103,0] = '07:59:28'
$arr[103,1] = '09:22'
$arr[103,2] = '16_SANTA ANA'
$arr[103,3] = 83
$arr[103,4] = 'LP1912'
$arr[104,0] = '07:47:32'
$arr[104,1] = '09:22'
$arr[104,2] = '17_ROMERO'
$arr[104,3] = 95
$arr[104,4] = 'LP1912'
$arr[105,0] = '08:57:13'
$arr[105,1] = '09:23'
$arr[105,2] = '16_SANTA ANA'
$arr[105,3] = 26
$arr[105,4] = 'LP1912'
$arr[106,0] = '07:47:32'
$arr[106,1] = '09:23'
$arr[106,2] = '11_ETCHEVERRY'
$arr[106,3] = 96
$arr[106,4] = 'LP1912'
$arr[107,0] = '08:57:13'
$arr[107,1] = '09:24'
$arr[107,2] = '11_ETCHEVERRY'
$arr[107,3] = 27
$arr[107,4] = 'LP1912'
$arr[108,0] = '08:21:50'
$arr[108,1] = '09:29'
$arr[108,2] = '16_SANTA ANA'
$arr[108,3] = 68
$arr[108,4] = 'LP1912'
$arr[109,0] = '07:47:32'
$arr[109,1] = '09:32'
$arr[109,2] = '15_ABASTO'
$arr[109,3] = 105
$arr[109,4] = 'LP1912'
$arr[110,0] = '07:47:32'
$arr[110,1] = '09:33'
$arr[110,2] = '10_OLMOS'
$arr[110,3] = 106
$arr[110,4] = 'LP1912'
$arr[111,0] = '08:39:44'
$arr[111,1] = '09:34'
$arr[111,2] = '23_HERNANDEZ'
$arr[111,3] = 55
$arr[111,4] = 'LP1912'
$arr[112,0] = '08:39:44'
$arr[112,1] = '09:34'
$arr[112,2] = '16_SANTA ANA'
$arr[112,3] = 55
$arr[112,4] = 'LP1912'
$arr[113,0] = '08:57:13'
$arr[113,1] = '09:35'
$arr[113,2] = '23_HERNANDEZ'
$arr[113,3] = 38
$arr[113,4] = 'LP1912'
$arr[114,0] = '08:57:13'
$arr[114,1] = '09:35'
$arr[114,2] = '16_SANTA ANA'
$arr[114,3] = 38
$arr[114,4] = 'LP1912'
$arr[115,0] = '08:21:50'
$arr[115,1] = '09:41'
$arr[115,2] = '215C_EL PATO'
$arr[115,3] = 80
$arr[115,4] = 'LP1912'
$arr[116,0] = '09:38:09'
$arr[116,1] = '09:41'
$arr[116,2] = '14_ABASTO'
$arr[116,3] = 3
$arr[116,4] = 'LP1912'
$arr[117,0] = '09:38:09'
$arr[117,1] = '09:41'
$arr[117,2] = '23_HERNANDEZ'
$arr[117,3] = 3
$arr[117,4] = 'LP1912'
$arr[118,0] = '07:47:32'
$arr[118,1] = '09:42'
$arr[118,2] = '215C_EL PATO'
$arr[118,3] = 115
$arr[118,4] = 'LP1912'
$arr[119,0] = '07:47:32'
$arr[119,1] = '09:43'
$arr[119,2] = '14_ABASTO'
$arr[119,3] = 116
$arr[119,4] = 'LP1912'
$arr[120,0] = '08:57:13'
$arr[120,1] = '09:44'
$arr[120,2] = '14_ABASTO'
$arr[120,3] = 47
$arr[120,4] = 'LP1912'
$arr[121,0] = '09:38:09'
$arr[121,1] = '09:47'
$arr[121,2] = '16_SANTA ANA'
$arr[121,3] = 9
$arr[121,4] = 'LP1912'
$arr[122,0] = '08:49:51'
$arr[122,1] = '09:52'
$arr[122,2] = '15_ABASTO'
$arr[122,3] = 63
$arr[122,4] = 'LP1912'
$arr[123,0] = '08:49:51'
$arr[123,1] = '09:53'
$arr[123,2] = '10_OLMOS'
$arr[123,3] = 64
$arr[123,4] = 'LP1912'
$arr[124,0] = '09:38:09'
$arr[124,1] = '09:59'
$arr[124,2] = '16_SANTA ANA'
$arr[124,3] = 21
$arr[124,4] = 'LP1912'
$arr[125,0] = '09:38:09'
$arr[125,1] = '10:04'
$arr[125,2] = '11_ETCHEVERRY'
$arr[125,3] = 26
$arr[125,4] = 'LP1912'
$arr[126,0] = '09:38:09'
$arr[126,1] = '10:05'
$arr[126,2] = '23_HERNANDEZ'
$arr[126,3] = 27
$arr[126,4] = 'LP1912'
$arr[127,0] = '08:39:44'
$arr[127,1] = '10:06'
$arr[127,2] = '10_OLMOS'
$arr[127,3] = 87
$arr[127,4] = 'LP1912'
$arr[128,0] = '08:21:50'
$arr[128,1] = '10:10'
$arr[128,2] = '16_P MOR-SANTA ANA'
$arr[128,3] = 109
$arr[128,4] = 'LP1912'
$arr[129,0] = '08:57:13'
$arr[129,1] = '10:11'
$arr[129,2] = '16_P MOR-SANTA ANA'
$arr[129,3] = 74
$arr[129,4] = 'LP1912'
$arr[130,0] = '08:21:50'
$arr[130,1] = '10:12'
$arr[130,2] = '15_ABASTO'
$arr[130,3] = 111
$arr[130,4] = 'LP1912'
$arr[131,0] = '09:38:09'
$arr[131,1] = '10:13'
$arr[131,2] = '10_OLMOS'
$arr[131,3] = 35
$arr[131,4] = 'LP1912'
$arr[132,0] = '08:49:51'
$arr[132,1] = '10:20'
$arr[132,2] = '26_HERNANDEZ'
$arr[132,3] = 91
$arr[132,4] = 'LP1912'
$arr[133,0] = '08:39:44'
$arr[133,1] = '10:21'
$arr[133,2] = '26_HERNANDEZ'
$arr[133,3] = 102
$arr[133,4] = 'LP1912'
$arr[134,0] = '08:39:44'
$arr[134,1] = '10:22'
$arr[134,2] = '17_ROMERO'
$arr[134,3] = 103
$arr[134,4] = 'LP1912'
$arr[135,0] = '09:38:09'
$arr[135,1] = '10:24'
$arr[135,2] = '11_ETCHEVERRY'
$arr[135,3] = 46
$arr[135,4] = 'LP1912'
$arr[136,0] = '08:39:44'
$arr[136,1] = '10:26'
$arr[136,2] = '215A_EL PATO'
$arr[136,3] = 107
$arr[136,4] = 'LP1912'
$arr[137,0] = '08:57:13'
$arr[137,1] = '10:27'
$arr[137,2] = '215A_EL PATO'
$arr[137,3] = 90
$arr[137,4] = 'LP1912'
$arr[138,0] = '10:26:41'
$arr[138,1] = '10:33'
$arr[138,2] = '10_OLMOS'
$arr[138,3] = 7
$arr[138,4] = 'LP1912'
$arr[139,0] = '10:26:41'
$arr[139,1] = '10:34'
$arr[139,2] = '23_HERNANDEZ'
$arr[139,3] = 8
$arr[139,4] = 'LP1912'
$arr[140,0] = '10:26:41'
$arr[140,1] = '10:34'
$arr[140,2] = '16_SANTA ANA'
$arr[140,3] = 8
$arr[140,4] = 'LP1912'
$arr[141,0] = '08:49:51'
$arr[141,1] = '10:41'
$arr[141,2] = '17_ROMERO'
$arr[141,3] = 112
$arr[141,4] = 'LP1912'
$arr[142,0] = '08:57:13'
$arr[142,1] = '10:42'
$arr[142,2] = '17_ROMERO'
$arr[142,3] = 105
$arr[142,4] = 'LP1912'
$arr[143,0] = '08:49:51'
$arr[143,1] = '10:43'
$arr[143,2] = '14_ABASTO'
$arr[143,3] = 114
$arr[143,4] = 'LP1912'
$arr[144,0] = '08:57:13'
$arr[144,1] = '10:44'
$arr[144,2] = '14_ABASTO'
$arr[144,3] = 107
$arr[144,4] = 'LP1912'
$arr[145,0] = '10:26:41'
$arr[145,1] = '10:46'
$arr[145,2] = '16_SANTA ANA'
$arr[145,3] = 20
$arr[145,4] = 'LP1912'
$arr[146,0] = '10:26:41'
$arr[146,1] = '10:52'
$arr[146,2] = '15_ABASTO'
$arr[146,3] = 26
$arr[146,4] = 'LP1912'
$arr[147,0] = '10:26:41'
$arr[147,1] = '10:53'
$arr[147,2] = '10_OLMOS'
$arr[147,3] = 27
$arr[147,4] = 'LP1912'
$arr[148,0] = '10:26:41'
$arr[148,1] = '10:56'
$arr[148,2] = '27_EL RETIRO'
$arr[148,3] = 30
$arr[148,4] = 'LP1912'
$arr[149,0] = '09:38:09'
$arr[149,1] = '10:58'
$arr[149,2] = '27_EL RETIRO'
$arr[149,3] = 80
$arr[149,4] = 'LP1912'
$arr[150,0] = '10:26:41'
$arr[150,1] = '11:01'
$arr[150,2] = '215C_EL PATO'
$arr[150,3] = 35
$arr[150,4] = 'LP1912'
$arr[151,0] = '09:38:09'
$arr[151,1] = '11:02'
$arr[151,2] = '215C_EL PATO'
$arr[151,3] = 84
$arr[151,4] = 'LP1912'
$arr[152,0] = '10:26:41'
$arr[152,1] = '11:03'
$arr[152,2] = '11_ETCHEVERRY'
$arr[152,3] = 37
$arr[152,4] = 'LP1912'
$arr[153,0] = '10:26:41'
$arr[153,1] = '11:04'
$arr[153,2] = '23_HERNANDEZ'
$arr[153,3] = 38
$arr[153,4] = 'LP1912'
$arr[154,0] = '10:26:41'
$arr[154,1] = '11:06'
$arr[154,2] = '16_P MOR-167 Y 521'
$arr[154,3] = 40
$arr[154,4] = 'LP1912'
$arr[155,0] = '09:38:09'
$arr[155,1] = '11:07'
$arr[155,2] = '16_P MOR-167 Y 521'
$arr[155,3] = 89
$arr[155,4] = 'LP1912'
$arr[156,0] = '10:26:41'
$arr[156,1] = '11:12'
$arr[156,2] = '15_ABASTO'
$arr[156,3] = 46
$arr[156,4] = 'LP1912'
$arr[157,0] = '10:26:41'
$arr[157,1] = '11:19'
$arr[157,2] = '86_EST CHICA-ESC AGRARIA'
$arr[157,3] = 53
$arr[157,4] = 'LP1912'
$arr[158,0] = '09:38:09'
$arr[158,1] = '11:20'
$arr[158,2] = '86_EST CHICA-ESC AGRARIA'
$arr[158,3] = 102
$arr[158,4] = 'LP1912'
$arr[159,0] = '09:38:09'
$arr[159,1] = '11:21'
$arr[159,2] = '26_HERNANDEZ'
$arr[159,3] = 103
$arr[159,4] = 'LP1912'
$arr[160,0] = '09:38:09'
$arr[160,1] = '11:27'
$arr[160,2] = '225_C ROCA-H SUR'
$arr[160,3] = 109
$arr[160,4] = 'LP1912'
$arr[161,0] = '09:38:09'
$arr[161,1] = '11:32'
$arr[161,2] = '81_EL PELIGRO'
$arr[161,3] = 114
$arr[161,4] = 'LP1912'
$arr[162,0] = '10:26:41'
$arr[162,1] = '11:35'
$arr[162,2] = '11_ETCHEVERRY'
$arr[162,3] = 69
$arr[162,4] = 'LP1912'
$arr[163,0] = '09:38:09'
$arr[163,1] = '11:36'
$arr[163,2] = '11_ETCHEVERRY'
$arr[163,3] = 118
$arr[163,4] = 'LP1912'
$arr[164,0] = '10:26:41'
$arr[164,1] = '11:41'
$arr[164,2] = '17_ROMERO'
$arr[164,3] = 75
$arr[164,4] = 'LP1912'
$arr[165,0] = '10:26:41'
$arr[165,1] = '11:51'
$arr[165,2] = '215B_EL PATO'
$arr[165,3] = 85
$arr[165,4] = 'LP1912'
$arr[166,0] = '10:26:41'
$arr[166,1] = '11:59'
$arr[166,2] = '225_GOMEZ'
$arr[166,3] = 93
$arr[166,4] = 'LP1912'
$arr[167,0] = '10:26:41'
$arr[167,1] = '12:02'
$arr[167,2] = '84_COLONIA URQUIZA-ESC 49'
$arr[167,3] = 96
$arr[167,4] = 'LP1912'
$arr[168,0] = '10:26:41'
$arr[168,1] = '12:06'
$arr[168,2] = '16_P MOR-SANTA ANA'
$arr[168,3] = 100
$arr[168,4] = 'LP1912'
$arr[169,0] = '10:26:41'
$arr[169,1] = '12:14'
$arr[169,2] = '17_ROMERO'
$arr[169,3] = 108
$arr[169,4] = 'LP1912'
$arr[170,0] = '10:26:41'
$arr[170,1] = '12:19'
$arr[170,2] = '14_ABASTO'
$arr[170,3] = 113
$arr[170,4] = 'LP1912'
$arr[171,0] = '10:26:41'
$arr[171,1] = '12:20'
$arr[171,2] = '215A_EL PATO'
$arr[171,3] = 114
$arr[171,4] = 'LP1912'
$arr[172,0] = '10:26:41'
$arr[172,1] = '12:21'
$arr[172,2] = '26_HERNANDEZ'
$arr[172,3] = 115
$arr[172,4] = 'LP1912'
$ws.Range("A6:E178").Value = $arr

# ----- Sheet: LP1912-215 -----
$ws = $wb.Worksheets.Item('LP1912-215')
$ws.Range("A2").Value = 'Última actualización: 10:26:41'
$ws.Range("A3").Value = 'Total filas: 23'
$arr = New-Object 'object[,]' 23,5
$arr[0,0] = '04:18:53'
$arr[0,1] = '04:45'
$arr[0,2] = '215A_EL PATO'
$arr[0,3] = 27
$arr[0,4] = 'LP1912'
$arr[1,0] = '04:40:33'
$arr[1,1] = '04:46'
$arr[1,2] = '215A_EL PATO'
$arr[1,3] = 6
$arr[1,4] = 'LP1912'
$arr[2,0] = '04:18:53'
$arr[2,1] = '05:34'
$arr[2,2] = '215B_EL PATO'
$arr[2,3] = 76
$arr[2,4] = 'LP1912'
$arr[3,0] = '04:54:03'
$arr[3,1] = '05:35'
$arr[3,2] = '215B_EL PATO'
$arr[3,3] = 41
$arr[3,4] = 'LP1912'
$arr[4,0] = '04:18:53'
$arr[4,1] = '06:11'
$arr[4,2] = '215A_EL PATO'
$arr[4,3] = 113
$arr[4,4] = 'LP1912'
$arr[5,0] = '04:54:03'
$arr[5,1] = '06:12'
$arr[5,2] = '215A_EL PATO'
$arr[5,3] = 78
$arr[5,4] = 'LP1912'
$arr[6,0] = '05:18:56'
$arr[6,1] = '06:46'
$arr[6,2] = '215C_EL PATO'
$arr[6,3] = 88
$arr[6,4] = 'LP1912'
$arr[7,0] = '04:54:03'
$arr[7,1] = '06:47'
$arr[7,2] = '215C_EL PATO'
$arr[7,3] = 113
$arr[7,4] = 'LP1912'
$arr[8,0] = '05:18:56'
$arr[8,1] = '07:11'
$arr[8,2] = '215A_EL PATO'
$arr[8,3] = 113
$arr[8,4] = 'LP1912'
$arr[9,0] = '06:15:04'
$arr[9,1] = '07:12'
$arr[9,2] = '215A_EL PATO'
$arr[9,3] = 57
$arr[9,4] = 'LP1912'
$arr[10,0] = '06:43:40'
$arr[10,1] = '07:51'
$arr[10,2] = '215D_EL PATO'
$arr[10,3] = 68
$arr[10,4] = 'LP1912'
$arr[11,0] = '06:15:04'
$arr[11,1] = '07:52'
$arr[11,2] = '215D_EL PATO'
$arr[11,3] = 97
$arr[11,4] = 'LP1912'
$arr[12,0] = '06:43:40'
$arr[12,1] = '08:23'
$arr[12,2] = '215B_EL PATO'
$arr[12,3] = 100
$arr[12,4] = 'LP1912'
$arr[13,0] = '07:20:40'
$arr[13,1] = '09:01'
$arr[13,2] = '215A_EL PATO'
$arr[13,3] = 101
$arr[13,4] = 'LP1912'
$arr[14,0] = '08:57:13'
$arr[14,1] = '09:02'
$arr[14,2] = '215A_EL PATO'
$arr[14,3] = 5
$arr[14,4] = 'LP1912'
$arr[15,0] = '08:21:50'
$arr[15,1] = '09:41'
$arr[15,2] = '215C_EL PATO'
$arr[15,3] = 80
$arr[15,4] = 'LP1912'
$arr[16,0] = '07:47:32'
$arr[16,1] = '09:42'
$arr[16,2] = '215C_EL PATO'
$arr[16,3] = 115
$arr[16,4] = 'LP1912'
$arr[17,0] = '08:39:44'
$arr[17,1] = '10:26'
$arr[17,2] = '215A_EL PATO'
$arr[17,3] = 107
$arr[17,4] = 'LP1912'
$arr[18,0] = '08:57:13'
$arr[18,1] = '10:27'
$arr[18,2] = '215A_EL PATO'
$arr[18,3] = 90
$arr[18,4] = 'LP1912'
$arr[19,0] = '10:26:41'
$arr[19,1] = '11:01'
$arr[19,2] = '215C_EL PATO'
$arr[19,3] = 35
$arr[19,4] = 'LP1912'
$arr[20,0] = '09:38:09'
$arr[20,1] = '11:02'
$arr[20,2] = '215C_EL PATO'
$arr[20,3] = 84
$arr[20,4] = 'LP1912'
$arr[21,0] = '10:26:41'
$arr[21,1] = '11:51'
$arr[21,2] = '215B_EL PATO'
$arr[21,3] = 85
$arr[21,4] = 'LP1912'
$arr[22,0] = '10:26:41'
$arr[22,1] = '12:20'
$arr[22,2] = '215A_EL PATO'
$arr[22,3] = 114
$arr[22,4] = 'LP1912'
$ws.Range("A6:E28").Value = $arr

# ----- Sheet: 6203-6173 -----
$ws = $wb.Worksheets.Item('6203-6173')
$ws.Range("A2").Value = 'Última actualización: 10:26:41'
$ws.Range("A3").Value = 'Total filas: 31'
$arr = New-Object 'object[,]' 31,5
$arr[0,0] = '04:18:53'
$arr[0,1] = '05:43'
$arr[0,2] = '215A_LA PLATA'
$arr[0,3] = 85
$arr[0,4] = 'L6173'
$arr[1,0] = '04:40:33'
$arr[1,1] = '05:44'
$arr[1,2] = '215A_LA PLATA'
$arr[1,3] = 64
$arr[1,4] = 'L6173'
$arr[2,0] = '04:18:53'
$arr[2,1] = '06:08'
$arr[2,2] = '215A_LA PLATA'
$arr[2,3] = 110
$arr[2,4] = 'L6173'
$arr[3,0] = '04:40:33'
$arr[3,1] = '06:09'
$arr[3,2] = '215A_LA PLATA'
$arr[3,3] = 89
$arr[3,4] = 'L6173'
$arr[4,0] = '05:18:56'
$arr[4,1] = '06:32'
$arr[4,2] = '215C_LA PLATA'
$arr[4,3] = 74
$arr[4,4] = 'L6203'
$arr[5,0] = '04:40:33'
$arr[5,1] = '06:33'
$arr[5,2] = '215C_LA PLATA'
$arr[5,3] = 113
$arr[5,4] = 'L6203'
$arr[6,0] = '05:18:56'
$arr[6,1] = '06:59'
$arr[6,2] = '215B_LP-P MOR-1 Y 57'
$arr[6,3] = 101
$arr[6,4] = 'L6173'
$arr[7,0] = '06:15:04'
$arr[7,1] = '07:00'
$arr[7,2] = '215B_LP-P MOR-1 Y 57'
$arr[7,3] = 45
$arr[7,4] = 'L6173'
$arr[8,0] = '06:57:30'
$arr[8,1] = '07:01'
$arr[8,2] = '215B_LP-P MOR-1 Y 57'
$arr[8,3] = 4
$arr[8,4] = 'L6173'
$arr[9,0] = '05:49:40'
$arr[9,1] = '07:07'
$arr[9,2] = '215B_LP-P MOR-1 Y 57'
$arr[9,3] = 78
$arr[9,4] = 'L6173'
$arr[10,0] = '05:49:40'
$arr[10,1] = '07:35'
$arr[10,2] = '215A_LA PLATA'
$arr[10,3] = 106
$arr[10,4] = 'L6173'
$arr[11,0] = '07:20:40'
$arr[11,1] = '07:37'
$arr[11,2] = '215A_LA PLATA'
$arr[11,3] = 17
$arr[11,4] = 'L6173'
$arr[12,0] = '06:43:40'
$arr[12,1] = '08:06'
$arr[12,2] = '215C_LA PLATA'
$arr[12,3] = 83
$arr[12,4] = 'L6203'
$arr[13,0] = '06:15:04'
$arr[13,1] = '08:07'
$arr[13,2] = '215C_LA PLATA'
$arr[13,3] = 112
$arr[13,4] = 'L6203'
$arr[14,0] = '07:20:40'
$arr[14,1] = '08:09'
$arr[14,2] = '215C_LA PLATA'
$arr[14,3] = 49
$arr[14,4] = 'L6203'
$arr[15,0] = '06:57:30'
$arr[15,1] = '08:10'
$arr[15,2] = '215C_LA PLATA'
$arr[15,3] = 73
$arr[15,4] = 'L6203'
$arr[16,0] = '07:47:32'
$arr[16,1] = '08:13'
$arr[16,2] = '215C_LA PLATA'
$arr[16,3] = 26
$arr[16,4] = 'L6203'
$arr[17,0] = '07:59:28'
$arr[17,1] = '08:17'
$arr[17,2] = '215C_LA PLATA'
$arr[17,3] = 18
$arr[17,4] = 'L6203'
$arr[18,0] = '06:57:30'
$arr[18,1] = '08:35'
$arr[18,2] = '215A_LA PLATA'
$arr[18,3] = 98
$arr[18,4] = 'L6173'
$arr[19,0] = '07:59:28'
$arr[19,1] = '08:36'
$arr[19,2] = '215A_LA PLATA'
$arr[19,3] = 37
$arr[19,4] = 'L6173'
$arr[20,0] = '06:43:40'
$arr[20,1] = '08:38'
$arr[20,2] = '215A_LA PLATA'
$arr[20,3] = 115
$arr[20,4] = 'L6173'
$arr[21,0] = '08:21:50'
$arr[21,1] = '08:42'
$arr[21,2] = '215A_LA PLATA'
$arr[21,3] = 21
$arr[21,4] = 'L6173'
$arr[22,0] = '08:39:44'
$arr[22,1] = '08:44'
$arr[22,2] = '215A_LA PLATA'
$arr[22,3] = 5
$arr[22,4] = 'L6173'
$arr[23,0] = '07:20:40'
$arr[23,1] = '09:08'
$arr[23,2] = '215D_LA PLATA'
$arr[23,3] = 108
$arr[23,4] = 'L6203'
$arr[24,0] = '07:47:32'
$arr[24,1] = '09:09'
$arr[24,2] = '215D_LA PLATA'
$arr[24,3] = 82
$arr[24,4] = 'L6203'
$arr[25,0] = '08:21:50'
$arr[25,1] = '10:02'
$arr[25,2] = '215B_LP-P MOR-40 Y 115'
$arr[25,3] = 101
$arr[25,4] = 'L6173'
$arr[26,0] = '08:57:13'
$arr[26,1] = '10:03'
$arr[26,2] = '215B_LP-P MOR-40 Y 115'
$arr[26,3] = 66
$arr[26,4] = 'L6173'
$arr[27,0] = '08:57:13'
$arr[27,1] = '10:54'
$arr[27,2] = '215A_LA PLATA'
$arr[27,3] = 117
$arr[27,4] = 'L6173'
$arr[28,0] = '10:26:41'
$arr[28,1] = '11:13'
$arr[28,2] = '215C_LA PLATA'
$arr[28,3] = 47
$arr[28,4] = 'L6203'
$arr[29,0] = '09:38:09'
$arr[29,1] = '11:14'
$arr[29,2] = '215C_LA PLATA'
$arr[29,3] = 96
$arr[29,4] = 'L6203'
$arr[30,0] = '10:26:41'
$arr[30,1] = '12:04'
$arr[30,2] = '215A_LA PLATA'
$arr[30,3] = 98
$arr[30,4] = 'L6173'
$ws.Range("A6:E36").Value = $arr

